# Apply the StructureDefinition-employee-division.xlsx update:
#  - Rebrand URL/Publisher from ibm.com/Alvearie to linuxforhealth.org/LinuxForHealth
#  - Bump Version 7.0.0 -> 8.0.0
#  - Update Date to the new publish timestamp
#  - Clear the stray Constraint(s) text that had been duplicated on the root
#    "Extension" row of the Elements sheet (it now only belongs to the
#    "Extension.extension" row, which already carries the correct text)

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-division"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
